$wb = $excel.ActiveWorkbook

# Helper: write a value as a genuine shared string (avoids Excel's
# auto-numeric-coercion for numeric-looking text like "0.544") by
# round-tripping through a text formula + paste-values.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# ------------------------------------------------------------------
# Add the two new worksheets after "bootstrapping"
# ------------------------------------------------------------------
$bootstrapping = $wb.Worksheets.Item("bootstrapping")

$wsRSquare = $wb.Worksheets.Add($null, $bootstrapping)
$wsRSquare.Name = "r square"

$wsBlindfold = $wb.Worksheets.Add($null, $wsRSquare)
$wsBlindfold.Name = "blindfold"

# ------------------------------------------------------------------
# "r square" sheet content
# ------------------------------------------------------------------
$wsRSquare.Range("B1").Value = "R Square"
$wsRSquare.Range("C1").Value = "R Square Adjusted"

$wsRSquare.Range("A2").Value = "PK (Y)"
Set-TextValue $wsRSquare.Range("B2") "0.544"
Set-TextValue $wsRSquare.Range("C2") "0.535"

# ------------------------------------------------------------------
# "blindfold" sheet content
# ------------------------------------------------------------------
$wsBlindfold.Range("B1").Value = "SSO"
$wsBlindfold.Range("C1").Value = "SSE"
$wsBlindfold.Range("D1").Value = "Q² (=1-SSE/SSO)"

$rows = @(
    @{ Row = 2; Label = "BK (X3)";       B = 2250000; C = 2250000; D = $null },
    @{ Row = 3; Label = "BK > D > PK";   B = 375000;  C = 375000;  D = $null },
    @{ Row = 4; Label = "D (Z)";         B = 4125000; C = 4125000; D = $null },
    @{ Row = 5; Label = "P (X1)";        B = 3750000; C = 3750000; D = $null },
    @{ Row = 6; Label = "P > D > PK";    B = 375000;  C = 375000;  D = $null },
    @{ Row = 7; Label = "PK (Y)";        B = 3000000; C = 1772020; D = "0.409" },
    @{ Row = 8; Label = "WB (X2)";       B = 2250000; C = 2250000; D = $null },
    @{ Row = 9; Label = "WB > D > PK";   B = 375000;  C = 375000;  D = $null }
)

foreach ($r in $rows) {
    $row = $r.Row
    $wsBlindfold.Range("A$row").Value = $r.Label

    $bCell = $wsBlindfold.Range("B$row")
    $bCell.NumberFormat = "#,##0"
    $bCell.Value = $r.B

    $cCell = $wsBlindfold.Range("C$row")
    $cCell.NumberFormat = "#,##0"
    $cCell.Value = $r.C

    if ($r.D) {
        Set-TextValue $wsBlindfold.Range("D$row") $r.D
    }
}

# ------------------------------------------------------------------
# Selections / active sheet
# ------------------------------------------------------------------
[void]$bootstrapping.Range("A1:G9").Select()
[void]$wsRSquare.Range("A1:D3").Select()
[void]$wsBlindfold.Range("A1:E10").Select()

$wsBlindfold.Activate()
